# Added filtering options for the Component Analysis
#
# The evaluation table is a rolling window of the latest 10 quarters
# (rows 2-11). A new quarter's evaluation data is now included, which
# shifts the existing rows down by one (row 2 -> row 3, row 3 -> row 4,
# etc.) and drops the oldest row (the previous row 11). A new row of
# values (for the newest observation, with the sample size "N" one
# larger than the previous top row) is then written into row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-11) down by one row, starting from the
# bottom so values are not overwritten before being read.
for ($row = 11; $row -ge 3; $row--) {
    $srcRow = $row - 1
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($row, 3).Value = $ws.Cells.Item($srcRow, 3).Value2
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($srcRow, 4).Value2
    $ws.Cells.Item($row, 5).Value = $ws.Cells.Item($srcRow, 5).Value2
    $ws.Cells.Item($row, 6).Value = $ws.Cells.Item($srcRow, 6).Value2
    $ws.Cells.Item($row, 7).Value = $ws.Cells.Item($srcRow, 7).Value2
}

# Write the newest observation's evaluation statistics into row 2.
$ws.Range("B2").Value = -0.01944167324622808
$ws.Range("C2").Value = 1.290629691791469
$ws.Range("D2").Value = 10.351458416846
$ws.Range("E2").Value = 3.21736824389842
$ws.Range("F2").Value = 3.252860963509644
$ws.Range("G2").Value = 46
